# final version of dvp
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    # Order matches first-use order of new shared strings in the target file:
    # sex(M) -> surface(wet) -> id/sub#(sub007) -> iDAPT#(iDAPT343) -> observor(Wu) -> N/A(reused)
    $ws.Cells.Item($r, 9).Value  = "M"             # I: sex
    $ws.Cells.Item($r, 12).Value = "wet"           # L: surface
    $ws.Cells.Item($r, 3).Value  = "sub007"        # C: id
    $ws.Cells.Item($r, 4).Value  = "sub007"        # D: sub#
    $ws.Cells.Item($r, 8).Value  = "iDAPT343"      # H: iDAPT#
    $ws.Cells.Item($r, 27).Value = "Wu"            # AA: observor

    $ws.Cells.Item($r, 5).Value  = "N/A"           # E: brand
    $ws.Cells.Item($r, 6).Value  = "N/A"           # F: name
    $ws.Cells.Item($r, 7).Value  = "N/A"           # G: Style#
    $ws.Cells.Item($r, 11).Value = "N/A"           # K: order
    $ws.Cells.Item($r, 13).Value = "N/A"           # M: repeats
    $ws.Cells.Item($r, 26).Value = "N/A"           # Z: compare
    $ws.Cells.Item($r, 28).Value = "N/A"           # AB: session

    $ws.Cells.Item($r, 10).Value = 12              # J: Footwear size

    $ws.Cells.Item($r, 15).Value = 2               # O: uphill
    $ws.Cells.Item($r, 16).Value = 2               # P: downhill
    $ws.Cells.Item($r, 17).Value = 3               # Q: first slip
    $ws.Cells.Item($r, 18).Value = 8               # R: pre slip
    $ws.Cells.Item($r, 19).Value = 9               # S: slip
    $ws.Cells.Item($r, 20).Value = 9               # T: thermal
    $ws.Cells.Item($r, 21).Value = 9               # U: fit
    $ws.Cells.Item($r, 22).Value = 9               # V: heaviness
    $ws.Cells.Item($r, 23).Value = 9               # W: overall
    $ws.Cells.Item($r, 24).Value = 9               # X: easy take off
    $ws.Cells.Item($r, 25).Value = 9               # Y: use

    $ws.Cells.Item($r, 29).Value = 43864                 # AC: date
    $ws.Cells.Item($r, 30).Value = 0.62708333333333333   # AD: time

    # AE (air temp), AF (ice temp), AG (RH) no longer present - clear them
    $ws.Cells.Item($r, 31).ClearContents()
    $ws.Cells.Item($r, 32).ClearContents()
    $ws.Cells.Item($r, 33).ClearContents()
}

$ws.Range("C3:AG3").Select()
